# Weekly update: insert a new price record at row 118 (Cebollín, Terminal La
# Palmera de La Serena), pushing the existing rows 118-227 down to 119-228.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 118 - this shifts rows 118:227 down to 119:228
# and naturally extends the sheet dimension to A1:R228.
$ws.Rows.Item(118).Insert()

# Populate the newly inserted row with this week's data point.
$ws.Cells.Item(118, 1).Value  = 8
$ws.Cells.Item(118, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(118, 3).Value  = "Coquimbo"
$ws.Cells.Item(118, 4).Value  = 44790
$ws.Cells.Item(118, 5).Value  = 4
$ws.Cells.Item(118, 6).Value  = 100112037
$ws.Cells.Item(118, 7).Value  = "Cebollín"
$ws.Cells.Item(118, 8).Value  = "Sin especificar"
$ws.Cells.Item(118, 9).Value  = "Primera"
$ws.Cells.Item(118, 10).Value = 1400
$ws.Cells.Item(118, 11).Value = 1400
$ws.Cells.Item(118, 12).Value = 1600
$ws.Cells.Item(118, 13).Value = 1500
$ws.Cells.Item(118, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(118, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(118, 16).Value = 250
$ws.Cells.Item(118, 17).Value = 6
$ws.Cells.Item(118, 18).Value = "Hortaliza"
